$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Add a new worksheet named "metadata" right after the existing "data" sheet
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"
$meta.Range("B1:G1").Font.Bold = $true
$meta.Range("B1:G1").HorizontalAlignment = -4108
$meta.Range("B1:G1").VerticalAlignment = -4160
$meta.Range("B1:G1").Borders.LineStyle = 1

# Data row
$meta.Range("A2").Value = 0
$meta.Range("A2").Font.Bold = $true
$meta.Range("A2").HorizontalAlignment = -4108
$meta.Range("A2").VerticalAlignment = -4160
$meta.Range("A2").Borders.LineStyle = 1
$meta.Range("B2").Value = "Hereditary angioedema"
$meta.Range("C2").Value = 226
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.1"
$meta.Range("E2").Value = "2021-07-14T12:10:25.438690Z"
$meta.Range("F2").Value = "2021-10-05 14:33:56.781787"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/226/?format=json"

# Update the time_taken timestamps on the "data" sheet
$data.Range("F2").Value = "2021-10-05 14:33:56.785509"
$data.Range("F3").Value = "2021-10-05 14:33:56.785517"
$data.Range("F4").Value = "2021-10-05 14:33:56.785520"
$data.Range("F5").Value = "2021-10-05 14:33:56.785523"
$data.Range("F6").Value = "2021-10-05 14:33:56.785526"
$data.Range("F7").Value = "2021-10-05 14:33:56.785528"
$data.Range("F8").Value = "2021-10-05 14:33:56.785531"

# Keep the "data" sheet as the active tab
$data.Activate()
